$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5; existing rows 5:90 shift down to 6:91.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with a new data record (same shape as the
# other rows in this table), carrying the product/market constants and the
# new date + volume values.
$ws.Cells.Item(5, 1).Value = 5
$ws.Cells.Item(5, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(5, 3).Value = "Maule"
$ws.Cells.Item(5, 4).Value = 44515
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(5, 6).Value = "Fruta"
$ws.Cells.Item(5, 7).Value = 100108
$ws.Cells.Item(5, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(5, 9).Value = 100108002
$ws.Cells.Item(5, 10).Value = "Mango"
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Primera"
$ws.Cells.Item(5, 13).Value = 210
$ws.Cells.Item(5, 14).Value = 7000
$ws.Cells.Item(5, 15).Value = 7000
$ws.Cells.Item(5, 16).Value = 7000
$ws.Cells.Item(5, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(5, 18).Value = "Perú"
$ws.Cells.Item(5, 19).Value = 1750
$ws.Cells.Item(5, 20).Value = 4

# Ensure the date cell keeps the date-formatted style used by the rest of
# column D (style index 2 in the original workbook).
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat
